$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep being treated as text, matching the
# original inline-string cell type instead of being auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '29.935.87'
$ws.Range('E2').Value = '  -1.02%  '
$ws.Range('D3').Value = '1.918.56'
$ws.Range('E3').Value = '  +1.19%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '320.03'
$ws.Range('E5').Value = '  -1.56%  '
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').Value = '0.5059'
$ws.Range('E7').Value = '  -2.27%  '
$ws.Range('D8').Value = '0.4033'
$ws.Range('E8').Value = '  +0.67%  '
$ws.Range('D9').Value = '0.08335'
$ws.Range('E9').Value = '  -0.71%  '
$ws.Range('D10').Value = '1.103'
$ws.Range('E10').Value = '  -1.00%  '
$ws.Range('B11').Value = 'Solana'
$ws.Range('C11').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D11').Value = '23.79'
$ws.Range('E11').Value = '  +2.43%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.910.34'
$ws.Range('E12').Value = '  +0.89%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '6.401'
$ws.Range('E13').Value = '  -0.62%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').Value = '7.233'
$ws.Range('E14').Value = '  -1.28%  '
$ws.Range('B15').Value = 'BinanceUSD'
$ws.Range('C15').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D15').Value = '1.005'
$ws.Range('E15').Value = '  +0.41%  '
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D16').Value = '92.27'
$ws.Range('E16').Value = '  -2.07%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').Value = '0.00001097'
$ws.Range('E17').Value = '  -1.12%  '
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').Value = '0.06509'
$ws.Range('E18').Value = '  -1.97%  '
$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D19').Value = '18.22'
$ws.Range('E19').Value = '  +0.03%  '
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').Value = '1.001'
$ws.Range('E20').Value = '  +0.13%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = '5.947'
$ws.Range('E21').Value = '  -0.08%  '
$ws.Range('B22').Value = 'WrappedBTC'
$ws.Range('C22').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D22').Value = '29.957.55'
$ws.Range('E22').Value = '  -0.90%  '
$ws.Range('B23').Value = 'Cosmos'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D23').Value = '11.32'
$ws.Range('E23').Value = '  +0.26%  '
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').Value = '2.195'
$ws.Range('E24').Value = '  -1.43%  '
$ws.Range('B25').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C25').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D25').Value = '2.143.30'
$ws.Range('E25').Value = '  +1.47%  '
$ws.Range('D26').Value = '22.09'
$ws.Range('E26').Value = '  +2.39%  '
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').Value = '162.42'
$ws.Range('E27').Value = '  +0.46%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').Value = '2.318'
$ws.Range('E28').Value = '  -1.60%  '
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').Value = '129.51'
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('B30').Value = 'ImmutableX'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D30').Value = '1.127'
$ws.Range('E30').Value = '  +3.20%  '
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').Value = '0.1038'
$ws.Range('E31').Value = '  -1.51%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '5.966'
$ws.Range('E32').Value = '  -1.99%  '
$ws.Range('B33').Value = 'HuobiToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D33').Value = '3.813'
$ws.Range('E33').Value = '  +1.75%  '
$ws.Range('B34').Value = 'VeChain'
$ws.Range('C34').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D34').Value = '0.02448'
$ws.Range('E34').Value = '  -1.83%  '
$ws.Range('B35').Value = 'InternetComputer(DFINITY)'
$ws.Range('C35').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D35').Value = '5.391'
$ws.Range('E35').Value = '  +2.19%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').Value = '0.06420'
$ws.Range('E36').Value = '  -1.87%  '
$ws.Range('B37').Value = 'Algorand'
$ws.Range('C37').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D37').Value = '0.2156'
$ws.Range('E37').Value = '  -1.77%  '
$ws.Range('B38').Value = 'ARBITRUM'
$ws.Range('C38').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D38').Value = '1.193'
$ws.Range('E38').Value = '  -2.37%  '
$ws.Range('B39').Value = 'TheSandbox'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D39').Value = '0.6478'
$ws.Range('E39').Value = '  -0.36%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = '8.691'
$ws.Range('E40').Value = '  -0.53%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').Value = '1.219'
$ws.Range('E41').Value = '  -0.91%  '
$ws.Range('D42').Value = '11.37'
$ws.Range('E42').Value = '  -3.66%  '
$ws.Range('B43').Value = 'NEARProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D43').Value = '2.215'
$ws.Range('E43').Value = '  +7.90%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = '13.44'
$ws.Range('E44').Value = '  +1.34%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').Value = '0.6040'
$ws.Range('E45').Value = '  -0.99%  '
$ws.Range('B46').Value = 'PancakeSwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D46').Value = '3.637'
$ws.Range('E46').Value = '  -1.75%  '
$ws.Range('B47').Value = 'EOS'
$ws.Range('C47').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D47').Value = '1.209'
$ws.Range('E47').Value = '  -2.19%  '
$ws.Range('D48').Value = '122.17'
$ws.Range('E48').Value = '  -1.88%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').Value = '78.83'
$ws.Range('E49').Value = '  -0.25%  '
$ws.Range('B50').Value = 'WEMIXTOKEN'
$ws.Range('C50').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D50').Value = '1.130'
$ws.Range('E50').Value = '  -2.90%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '0.06798'
$ws.Range('E51').Value = '  -0.89%  '
